$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the "Execution" value for the two data rows to "Manual"
$ws.Range("E2").Value = "Manual"
$ws.Range("E3").Value = "Manual"

# Match the style used by the other data cells in the row (wrap text style)
$ws.Range("E2").WrapText = $true
$ws.Range("E3").WrapText = $true

# Update the active selection to E4, as reflected in the saved view state
$ws.Range("E4").Select()
